$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.677.18"
$ws.Range("E2").Value = "  +4.35%  "

# Row 3
$ws.Range("D3").Value = "1.750.95"
$ws.Range("E3").Value = "  +4.76%  "

# Row 4
$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "247.30"
$ws.Range("E5").Value = "  +3.45%  "

# Row 6
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("D7").Value = "0.4809"
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "0.2715"
$ws.Range("E8").Value = "  +3.11%  "

# Row 9
$ws.Range("D9").Value = "0.06261"
$ws.Range("E9").Value = "  +1.23%  "

# Row 10
$ws.Range("D10").Value = "1.746.65"
$ws.Range("E10").Value = "  +4.47%  "

# Row 11
$ws.Range("D11").Value = "0.07121"
$ws.Range("E11").Value = "  +1.49%  "

# Row 12
$ws.Range("D12").Value = "15.88"
$ws.Range("E12").Value = "  +6.50%  "

# Row 13
$ws.Range("D13").Value = "0.6231"
$ws.Range("E13").Value = "  +5.05%  "

# Row 14
$ws.Range("D14").Value = "4.519"
$ws.Range("E14").Value = "  +2.90%  "

# Row 15
$ws.Range("D15").Value = "77.37"
$ws.Range("E15").Value = "  +2.63%  "

# Row 16
$ws.Range("D16").Value = "0.9985"
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").Value = "26.668.35"
$ws.Range("E17").Value = "  +4.33%  "

# Row 18
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  -0.06%  "

# Row 19
$ws.Range("D19").Value = "0.000006904"
$ws.Range("E19").Value = "  +1.72%  "

# Row 20
$ws.Range("D20").Value = "11.76"
$ws.Range("E20").Value = "  +2.39%  "

# Row 21
$ws.Range("D21").Value = "1.969.28"
$ws.Range("E21").Value = "  +4.58%  "

# Row 22 (only D changes)
$ws.Range("D22").Value = "4.659"

# Row 23
$ws.Range("D23").Value = "8.887"
$ws.Range("E23").Value = "  +1.65%  "

# Row 24
$ws.Range("D24").Value = "5.358"
$ws.Range("E24").Value = "  +1.26%  "

# Row 25
$ws.Range("D25").Value = "136.42"
$ws.Range("E25").Value = "  -0.38%  "

# Row 26
$ws.Range("D26").Value = "15.53"
$ws.Range("E26").Value = "  +2.99%  "

# Row 27
$ws.Range("D27").Value = "1.839"
$ws.Range("E27").Value = "  +6.20%  "

# Row 28
$ws.Range("D28").Value = "1.407"
$ws.Range("E28").Value = "  +1.32%  "

# Row 29
$ws.Range("D29").Value = "107.95"
$ws.Range("E29").Value = "  +2.89%  "

# Row 30
$ws.Range("D30").Value = "4.032"
$ws.Range("E30").Value = "  +0.95%  "

# Row 31
$ws.Range("D31").Value = "3.775"
$ws.Range("E31").Value = "  +3.23%  "

# Row 32
$ws.Range("D32").Value = "0.07903"
$ws.Range("E32").Value = "  +1.06%  "

# Row 33
$ws.Range("D33").Value = "0.04592"
$ws.Range("E33").Value = "  +8.43%  "

# Row 34 (only E changes)
$ws.Range("E34").Value = "  -0.28%  "

# Row 35
$ws.Range("D35").Value = "0.6402"
$ws.Range("E35").Value = "  +4.87%  "

# Row 36
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +5.22%  "

# Row 37
$ws.Range("D37").Value = "0.9494"
$ws.Range("E37").Value = "  +10.65%  "

# Row 38
$ws.Range("D38").Value = "114.29"
$ws.Range("E38").Value = "  +18.81%  "

# Row 39
$ws.Range("D39").Value = "2.504"
$ws.Range("E39").Value = "  -3.71%  "

# Row 40
$ws.Range("D40").Value = "1.998"
$ws.Range("E40").Value = "  +6.66%  "

# Row 41
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.21%  "

# Row 42 and 43 swap contents (VeChain <-> FraxShare), with new price/volume values
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.766"
$ws.Range("E42").Value = "  +18.88%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.01516"
$ws.Range("E43").Value = "  +2.79%  "

# Row 44
$ws.Range("D44").Value = "0.3933"
$ws.Range("E44").Value = "  +3.94%  "

# Row 45
$ws.Range("D45").Value = "6.735"
$ws.Range("E45").Value = "  +8.09%  "

# Row 46
$ws.Range("D46").Value = "0.1207"
$ws.Range("E46").Value = "  +7.89%  "

# Row 47
$ws.Range("D47").Value = "0.05339"
$ws.Range("E47").Value = "  +1.53%  "

# Row 48
$ws.Range("D48").Value = "8.039"
$ws.Range("E48").Value = "  +8.84%  "

# Row 49
$ws.Range("D49").Value = "30.96"
$ws.Range("E49").Value = "  +3.68%  "

# Row 50
$ws.Range("D50").Value = "0.3468"
$ws.Range("E50").Value = "  +3.56%  "

# Row 51
$ws.Range("D51").Value = "51.92"
$ws.Range("E51").Value = "  +3.52%  "
